$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "dSF" (column F) values per re-pulled / re-calculated data.
# Row 9 (F9) is intentionally left unchanged (stays 0).
$values = @{
    2  = -4
    3  = -4
    4  = -3
    5  = -4
    6  = -4
    7  = -6
    8  = -5
    10 = -3
    11 = -2
    12 = -1
    13 = 3
    14 = -2
    15 = -2
    16 = 4
    17 = -2
    18 = 3
    19 = 4
    20 = 4
    21 = -6
    22 = -1
}

foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row]
}
